# Generate Report for Handback
#
# The handback CI run regenerated this status report: the two source
# files got new GUID-based names (and the zh-cn/de-de handoff/handback
# packages got a new content hash + refreshed timestamps). Every place
# in the workbook that spelled out the old file names / timestamps -
# plain cell text as well as each hyperlink's display text - needs to
# be swapped for the new value, while the hyperlink targets (URLs)
# themselves are left alone.

$oldGuid1 = "29e79f51-6ede-4853-a79d-4cea48aefdf7"
$oldGuid2 = "69bc8315-b512-49ab-a3b9-5d471a9f1a0e"
$newGuid1 = "1f641bcb-6eb4-4a1c-8351-6d8c6ffc8848"
$newGuid2 = "ffff1f6804c2-9b89-4dcb-86c9-d98c33ec6e5e"

$oldHash1 = "33aa57355ad5153a5270f755ce14a331139f09e6"
$oldHash2 = "3c07634513735dd3901ca696129c0c39357ebf3d"
$newHash  = "08af7258856ed5122eb31db14c922250c28858c4"

# Exact old-text -> new-text replacements (cell text and hyperlink
# display text both use these same literal strings).
$map = @{}
$map[$oldGuid1 + ".md"] = $newGuid1 + ".md"
$map[$oldGuid2 + ".md"] = $newGuid2 + ".md"

$map[$oldGuid1 + "." + $oldHash1 + ".zh-cn.xlf"] = $newGuid1 + "." + $newHash + ".zh-cn.xlf"
$map[$oldGuid2 + "." + $oldHash2 + ".zh-cn.xlf"] = $newGuid1 + "." + $newHash + ".zh-cn.xlf"

$map[$oldGuid1 + "." + $oldHash1 + ".de-de.xlf"] = $newGuid1 + "." + $newHash + ".de-de.xlf"
$map[$oldGuid2 + "." + $oldHash2 + ".de-de.xlf"] = $newGuid1 + "." + $newHash + ".de-de.xlf"

$map["2016-03-18 08:46:03"] = "2016-03-18 08:47:16"
$map["2016-03-18 08:46:22"] = "2016-03-18 08:47:32"
$map["2016-03-18 08:46:06"] = "2016-03-18 08:47:19"
$map["2016-03-18 08:46:30"] = "2016-03-18 08:47:37"

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {

    # Update plain cell text that exactly matches one of the old strings.
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        $text = $cell.Text
        if ($text -and $map.ContainsKey($text)) {
            $cell.Value = $map[$text]
        }
    }

    # Update the display text shown for each hyperlink (the hyperlink's
    # target URL / relationship is untouched).
    foreach ($hl in $ws.Hyperlinks) {
        $disp = $hl.TextToDisplay
        if ($disp -and $map.ContainsKey($disp)) {
            $hl.TextToDisplay = $map[$disp]
        }
    }
}
